$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'74.281.73"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +7.91%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.630.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +7.85%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'185.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +14.57%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'582.97"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +4.38%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  -0.18%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  +4.63%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  +19.48%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'2.628.54"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +7.82%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  +0.46%  "
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = "'  +8.10%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'4.76"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +3.51%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.0000190"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +7.15%  "
$ws.Range("E14").ClearFormats()
$ws.Range("B15").Value = "'WrappedBTC"
$ws.Range("B15").ClearFormats()
$ws.Range("C15").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C15").ClearFormats()
$ws.Range("D15").Value = "'74.117.43"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +7.85%  "
$ws.Range("E15").ClearFormats()
$ws.Range("B16").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B16").ClearFormats()
$ws.Range("C16").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C16").ClearFormats()
$ws.Range("D16").Value = "'3.108.99"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +7.65%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'26.28"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +13.01%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.628.84"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +7.76%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'9.10"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +31.01%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'11.84"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +11.70%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'372.05"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +9.70%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'2.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +18.89%  "
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'  +6.78%  "
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'  +0.21%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'70.08"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +4.48%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  +11.79%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'9.38"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +14.11%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'2.765.34"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +7.72%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +3.83%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'0.0₃0948"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +15.62%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'528.65"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +23.01%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +19.98%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'7.69"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +7.83%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "'  +8.97%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -0.11%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'163.10"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +1.94%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  +12.88%  "
$ws.Range("E37").ClearFormats()
$ws.Range("E38").Value = "'  +6.47%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'19.26"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +1.45%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  +0.03%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  +13.15%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.329"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +9.96%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'1.68"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +11.09%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'160.97"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +23.82%  "
$ws.Range("E44").ClearFormats()
$ws.Range("B45").Value = "'ImmutableX"
$ws.Range("B45").ClearFormats()
$ws.Range("C45").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C45").ClearFormats()
$ws.Range("D45").Value = "'1.19"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +11.90%  "
$ws.Range("E45").ClearFormats()
$ws.Range("B46").Value = "'dogwifhat"
$ws.Range("B46").ClearFormats()
$ws.Range("C46").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C46").ClearFormats()
$ws.Range("D46").Value = "'2.37"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +15.29%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'38.97"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +3.82%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.0856"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +18.76%  "
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'  +9.02%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.529"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +9.89%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'20.91"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +23.67%  "
$ws.Range("E51").ClearFormats()

Write-Host "Applied 92 cell updates"
